# Updates cryptos list values (prices and hourly volume % changes) to match
# the latest GitHub Actions scrape. Also reorders three rows (34-36) whose
# relative ranking changed: Binance-PegBSC-USD, Aptos, RenzoRestakedETH.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D/E store plain text (prices/percentages use dots as thousand
# separators, e.g. "67.779.53", which Excel would otherwise auto-parse as
# a number and mangle via float rounding). Mark the handful of new price
# values that look like ordinary decimals as Text before assigning them so
# they round-trip as the literal strings from the source feed.
$numericLookingCells = @("D4", "D5", "D6", "D13", "D20", "D21", "D24", "D28", "D31", "D33", "D34", "D35", "D38", "D42", "D44", "D46", "D47", "D49", "D51")
foreach ($cellref in $numericLookingCells) {
    $ws.Range($cellref).NumberFormat = "@"
}

$ws.Range("D2").Value = "67.779.53"
$ws.Range("E2").Value = "  +0.24%  "
$ws.Range("D3").Value = "3.811.60"
$ws.Range("E3").Value = "  +0.72%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "607.94"
$ws.Range("E5").Value = "  +2.10%  "
$ws.Range("D6").Value = "166.99"
$ws.Range("E6").Value = "  +0.23%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  +0.52%  "
$ws.Range("E9").Value = "  +0.77%  "
$ws.Range("E10").Value = "  -0.78%  "
$ws.Range("E11").Value = "  +0.59%  "
$ws.Range("E12").Value = "  -0.98%  "
$ws.Range("D13").Value = "36.07"
$ws.Range("E13").Value = "  -0.89%  "
$ws.Range("D14").Value = "4.448.20"
$ws.Range("E14").Value = "  +0.57%  "
$ws.Range("D15").Value = "3.801.74"
$ws.Range("E15").Value = "  +0.13%  "
$ws.Range("E16").Value = "  +0.51%  "
$ws.Range("D17").Value = "67.799.59"
$ws.Range("E17").Value = "  +0.33%  "
$ws.Range("E18").Value = "  +1.49%  "
$ws.Range("E19").Value = "  +0.46%  "
$ws.Range("D20").Value = "462.56"
$ws.Range("E20").Value = "  +1.37%  "
$ws.Range("D21").Value = "9.92"
$ws.Range("E21").Value = "  -2.69%  "
$ws.Range("E22").Value = "  +0.86%  "
$ws.Range("E23").Value = "  -2.42%  "
$ws.Range("D24").Value = "83.45"
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("E25").Value = "  +1.52%  "
$ws.Range("E26").Value = "  -1.17%  "
$ws.Range("E27").Value = "  +0.22%  "
$ws.Range("D28").Value = "10.05"
$ws.Range("E28").Value = "  -0.27%  "
$ws.Range("D29").Value = "3.959.26"
$ws.Range("E29").Value = "  +0.66%  "
$ws.Range("E30").Value = "  +0.59%  "
$ws.Range("D31").Value = "7.43"
$ws.Range("E31").Value = "  +1.93%  "
$ws.Range("E32").Value = "  +1.79%  "
$ws.Range("D33").Value = "29.63"
$ws.Range("E33").Value = "  -0.61%  "
$ws.Range("B34").Value = "Binance-PegBSC-USD"
$ws.Range("C34").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  +0.29%  "
$ws.Range("B35").Value = "Aptos"
$ws.Range("C35").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D35").Value = "9.09"
$ws.Range("E35").Value = "  -1.52%  "
$ws.Range("B36").Value = "RenzoRestakedETH"
$ws.Range("C36").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D36").Value = "3.750.61"
$ws.Range("E36").Value = "  +0.36%  "
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("D38").Value = "3.39"
$ws.Range("E38").Value = "  +1.75%  "
$ws.Range("E39").Value = "  -0.05%  "
$ws.Range("E40").Value = "  +0.26%  "
$ws.Range("E41").Value = "  +0.77%  "
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("D44").Value = "48.16"
$ws.Range("E44").Value = "  +2.32%  "
$ws.Range("E45").Value = "  +0.55%  "
$ws.Range("D46").Value = "43.19"
$ws.Range("E46").Value = "  -4.57%  "
$ws.Range("D47").Value = "28.19"
$ws.Range("E47").Value = "  +10.36%  "
$ws.Range("E48").Value = "  +0.17%  "
$ws.Range("D49").Value = "148.82"
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("E50").Value = "  +10.36%  "
$ws.Range("D51").Value = "1.85"
$ws.Range("E51").Value = "  +0.36%  "

# Re-normalise style on the cells we forced to Text above so they keep
# the same (default) cell style as every other data cell in the sheet —
# only the text content should differ from the original workbook.
$defaultStyle = $ws.Range("D7").Style
foreach ($cellref in $numericLookingCells) {
    $ws.Range($cellref).Style = $defaultStyle
}
